# Update the header row (row 1) column names to lowercase, keeping all
# other data untouched. This mirrors the commit:
# "Lo mismo que antes, edite los nombres de las columnas quitando las
#  mayusculas unicamente."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "filo_o_division"
$ws.Range("A1").Value = "reino"
$ws.Range("C1").Value = "clase"
$ws.Range("D1").Value = "orden"
$ws.Range("E1").Value = "familia"
$ws.Range("F1").Value = "genero"
$ws.Range("G1").Value = "epiteto_especifico"
$ws.Range("H1").Value = "nombre_comun"

# Match the selection/active-cell change seen in the diff.
$ws.Range("I5").Select()
